# Apply the data dictionary update: change the date format values from
# MM/DD/YYYY to YYYY-MM-DD for the three date fields (cdc_report_dt,
# onset_dt, pos_spec_dt) in rows 4, 9, and 10 of the "Values" column (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Dictionary")

$ws.Range("D4").Value = "YYYY-MM-DD"
$ws.Range("D9").Value = "YYYY-MM-DD"
$ws.Range("D10").Value = "YYYY-MM-DD"

# The two re-typed cells (D9/D10) picked up a "no fill" cell style on
# re-save (matches the saved file's xf index for these cells).
$ws.Range("D9").Interior.Pattern = -4142
$ws.Range("D10").Interior.Pattern = -4142

# Reflect the view state change captured in the saved file: selection
# landed on D7 after the edits.
$ws.Activate()
$ws.Range("D7").Select()

$wb.Save()
